$wb = $excel.ActiveWorkbook

# =========================================================================
# Sheet: 土地 (Land) -- add legislator/source metadata columns (I:O) and
# clean up existing text values, per issue #5 "property land done".
# =========================================================================
$ws = $wb.Worksheets.Item("土地")

# Re-write existing header row (B1:H1) with the new canonical column names
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "area"
$ws.Range("D1").Value = "share_portion"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "acquire_value"

# New metadata header columns
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

$hdr1 = $ws.Range("I1:O1")
$hdr1.Borders.LineStyle = 1
$hdr1.Font.Bold = $true
$hdr1.HorizontalAlignment = -4108
$hdr1.VerticalAlignment = -4160

# Clean up stray whitespace / punctuation in the existing data row
$ws.Range("B2").Value = "嘉義市嘉義市長竹段00060002地號"
$ws.Range("F2").Value = "96年02月09日"

# New metadata values for the data row
$ws.Range("I2").Value = "land"
$ws.Range("J2").Value = "normal"
$ws.Range("K2").Value = "'2013-12-02"
$ws.Range("L2").Value = "李俊俋"
$ws.Range("M2").Value = 1738
$ws.Range("N2").Value = "tmp52b51"
$ws.Range("O2").Value = 15

# =========================================================================
# Sheet: 建物 (Building) -- text cleanup only
# =========================================================================
$ws = $wb.Worksheets.Item("建物")
$ws.Range("B2").Value = "嘉義市嘉義市長竹段01946000建號"
$ws.Range("F2").Value = "96年02月09日"

# =========================================================================
# Sheet: 汽車 (Car) -- text cleanup only
# =========================================================================
$ws = $wb.Worksheets.Item("汽車")
$ws.Range("B2").Value = "MAZDA3"
$ws.Range("E2").Value = "98年09月16日"
$ws.Range("B3").Value = "MAZDAMPV"
$ws.Range("C3").Value = "'3000"
$ws.Range("E3").Value = "94年08月01日"
$ws.Range("E4").Value = "100年09月07日"

# =========================================================================
# Sheet: 存款 (Deposit) -- text cleanup only
# =========================================================================
$ws = $wb.Worksheets.Item("存款")
$ws.Range("B1").Value = "存放機構(應敘明分支機構）"
$ws.Range("B2").Value = "中華郵政股份有限公司台大郵局"
$ws.Range("B5").Value = "中華郵政股份有限公司嘉義中山路郵局"
$ws.Range("B10").Value = "中華郵政股份有限公司嘉義中山路郵局"
$ws.Range("B11").Value = "中華郵政股份有限公司嘉義中山路郵局"

# =========================================================================
# Sheet: 保險 (Insurance) -- text cleanup only
# =========================================================================
$ws = $wb.Worksheets.Item("保險")
$ws.Range("C2").Value = "富邦人壽心得意利率變動型年金保險"
$ws.Range("C3").Value = "月月金喜利率變動型養老保險"

# =========================================================================
# Sheet: 債務 (Debt) -- text cleanup only
# =========================================================================
$ws = $wb.Worksheets.Item("債務")
$ws.Range("D2").Value = "第一商業銀行嘉義市嘉義市東區文心街101巷8號"
$ws.Range("F2").Value = "96年02月12日"
$ws.Range("D3").Value = "玉山商業銀行嘉義市嘉義市東區文心街101巷8號"
$ws.Range("E3").Value = "'3893589"
$ws.Range("F3").Value = "96年02月12日"
